# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number: force text format first so
# Excel keeps them as text (matching the original inlineStr cells) instead of
# silently converting them to numeric values (which would also drop things like
# trailing zeros, e.g. "1.00" -> 1).
$textCells = @("D5", "D6", "D7", "D10", "D11", "D14", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D50")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.541.23'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '2.681.58'
$ws.Range("E3").Value = '  +3.01%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '582.94'
$ws.Range("E5").Value = '  +1.34%  '
$ws.Range("D6").Value = '145.47'
$ws.Range("E6").Value = '  +1.74%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.87%  '
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("D10").Value = '0.111'
$ws.Range("E10").Value = '  +4.78%  '
$ws.Range("D11").Value = '0.383'
$ws.Range("E11").Value = '  +4.01%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = '3.149.41'
$ws.Range("E13").Value = '  +2.95%  '
$ws.Range("D14").Value = '26.12'
$ws.Range("E14").Value = '  +7.09%  '
$ws.Range("D15").Value = '61.473.07'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("E16").Value = '  +4.55%  '
$ws.Range("D17").Value = '2.674.14'
$ws.Range("E17").Value = '  +2.68%  '
$ws.Range("D18").Value = '11.74'
$ws.Range("E18").Value = '  +2.79%  '
$ws.Range("D19").Value = '4.81'
$ws.Range("E19").Value = '  +4.20%  '
$ws.Range("D20").Value = '356.62'
$ws.Range("E20").Value = '  +2.85%  '
$ws.Range("D21").Value = '6.96'
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").Value = '0.525'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '64.79'
$ws.Range("E24").Value = '  +2.66%  '
$ws.Range("E25").Value = '  +3.26%  '
$ws.Range("D26").Value = '8.52'
$ws.Range("E26").Value = '  +5.96%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  +8.56%  '
$ws.Range("D29").Value = '0.0₃0832'
$ws.Range("E29").Value = '  +4.86%  '
$ws.Range("D30").Value = '6.98'
$ws.Range("E30").Value = '  +9.34%  '
$ws.Range("D31").Value = '169.74'
$ws.Range("E31").Value = '  +2.70%  '
$ws.Range("D32").Value = '0.998'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '20.23'
$ws.Range("E33").Value = '  +4.12%  '
$ws.Range("D34").Value = '1.15'
$ws.Range("E34").Value = '  +16.50%  '
$ws.Range("D35").Value = '4.72'
$ws.Range("E35").Value = '  +10.47%  '
$ws.Range("E36").Value = '  +6.43%  '
$ws.Range("D37").Value = '1.01'
$ws.Range("E37").Value = '  +20.09%  '
$ws.Range("D38").Value = '1.73'
$ws.Range("E38").Value = '  +6.37%  '
$ws.Range("D39").Value = '349.13'
$ws.Range("E39").Value = '  +11.34%  '
$ws.Range("E40").Value = '  +6.49%  '
$ws.Range("D41").Value = '38.63'
$ws.Range("E41").Value = '  +1.54%  '
$ws.Range("D42").Value = '5.41'
$ws.Range("E42").Value = '  +8.09%  '
$ws.Range("D43").Value = '21.05'
$ws.Range("E43").Value = '  +6.40%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").Value = '0.0583'
$ws.Range("E44").Value = '  +5.81%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '21.35'
$ws.Range("E45").Value = '  +6.73%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '136.67'
$ws.Range("E46").Value = '  +1.12%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.631'
$ws.Range("E47").Value = '  +4.28%  '
$ws.Range("E48").Value = '  +5.59%  '
$ws.Range("E49").Value = '  +1.07%  '
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").Value = '2.115.82'
$ws.Range("E51").Value = '  +4.77%  '
